$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text
$ws.Range("A1").Value = "The PC Target"
$ws.Range("B1").Value = "Responsible ICD Officer"
$ws.Range("C1").Value = "Status of The PC Target"

# Update the view/selection: top-left visible cell B1, active cell D1
$ws.Range("D1").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
